$wb = $excel.ActiveWorkbook

# The same edits apply to both the "展览" and "全部类型" worksheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: G2 ticket price becomes "已停售" (sold out / stopped selling) text
    $ws.Range("G2").Value = "已停售"

    # Row 3: F3 interest count 1459 -> 1461
    $ws.Range("F3").Value = 1461

    # Row 4: F4 interest count 954 -> 955
    $ws.Range("F4").Value = 955

    # Row 6: F6 interest count 2151 -> 2167
    $ws.Range("F6").Value = 2167

    # Row 7: F7 interest count 35 -> 36
    $ws.Range("F7").Value = 36

    # Row 8: F8 interest count 1309 -> 1324
    $ws.Range("F8").Value = 1324

    # Row 10: F10 interest count 129 -> 133
    $ws.Range("F10").Value = 133

    # Row 11: F11 interest count 41 -> 42
    $ws.Range("F11").Value = 42

    # Row 12: F12 interest count 315 -> 318
    $ws.Range("F12").Value = 318
}
